$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = -1.076622977925301
$ws.Range("J5").Value = 0.448376582304099
$ws.Range("K5").Value = 0.05765380539135707
$ws.Range("L5").Value = 2.483069796822584
